# Auto-generated: applies the 2023-09-11 cryptos-list refresh to Sheet1.
# Most cells are plain text assignments. A handful of Price (column D)
# values parse as valid numbers (e.g. "207.72"), so Excel would silently
# coerce them to the Number type on assignment; those go through a
# NumberFormat "@" (Text) round-trip first so they stay text, matching
# the original inline-string cells, then the cell style is reset back to
# "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.371.66"
$ws.Range("E2").Value = "  -2.43%  "
$ws.Range("D3").Value = "1.574.49"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("E4").Value = "  +0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "207.72"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -3.09%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  -4.45%  "
$ws.Range("E8").Value = "  -1.79%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.0609"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.52%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "17.82"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D12").Value = "1.792.85"
$ws.Range("E12").Value = "  -3.80%  "
$ws.Range("D13").Value = "1.575.51"
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("E14").Value = "  -3.22%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.507"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").Value = "25.367.55"
$ws.Range("E16").Value = "  -2.39%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "59.94"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").Value = "0.0₃0708"
$ws.Range("E18").Value = "  -4.02%  "
$ws.Range("E19").Value = "  +0.16%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "186.06"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("E21").Value = "  -1.97%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "9.31"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.28%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.90"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("E24").Value = "  +0.03%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "141.29"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("E26").Value = "  -2.83%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.70"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -4.82%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "14.92"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("E29").Value = "  -3.87%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.16"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.66%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0464"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.18%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.06"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.98%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.01"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.30%  "
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "1.084.31"
$ws.Range("E36").Value = "  -3.59%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.52%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.34"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("E39").Value = "  -2.30%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.779"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -9.07%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.494"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.30%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "94.44"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.69%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.07"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.83%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.707.26"
$ws.Range("E44").Value = "  -3.75%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.726"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -6.17%  "
$ws.Range("E46").Value = "  -4.86%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "52.76"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.62%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0507"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.408"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.41"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("E51").Value = "  -0.10%  "
